$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin/Link/Price/Volume(1h) refresh pulled from the latest coinranking.com snapshot.
# Rows 29-30, 34-35 and 39-51 also shuffle position as coins swap rank order.
$updates = @(
    @{ Row = 2; Price = "59.526.98"; Volume = "  -0.60%  " }
    @{ Row = 3; Price = "2.521.49"; Volume = "  -0.61%  " }
    @{ Row = 4; Price = "0.999"; Volume = "  -0.07%  " }
    @{ Row = 5; Price = "542.22"; Volume = "  -0.38%  " }
    @{ Row = 6; Price = "140.31"; Volume = "  -3.74%  " }
    @{ Row = 7; Price = "0.997"; Volume = "  +0.32%  " }
    @{ Row = 8; Price = "0.566"; Volume = "  -1.34%  " }
    @{ Row = 9; Price = "2.528.54"; Volume = "  -1.52%  " }
    @{ Row = 10; Price = "0.102"; Volume = "  +0.51%  " }
    @{ Row = 11; Volume = "  -0.14%  " }
    @{ Row = 12; Price = "5.40"; Volume = "  -3.41%  " }
    @{ Row = 13; Price = "0.353"; Volume = "  -2.57%  " }
    @{ Row = 14; Price = "2.966.01"; Volume = "  -0.56%  " }
    @{ Row = 15; Price = "23.39"; Volume = "  -0.98%  " }
    @{ Row = 16; Price = "59.400.88"; Volume = "  -0.62%  " }
    @{ Row = 17; Price = "0.0000142"; Volume = "  -1.51%  " }
    @{ Row = 18; Price = "2.512.95"; Volume = "  -1.78%  " }
    @{ Row = 19; Price = "11.13"; Volume = "  -1.28%  " }
    @{ Row = 20; Price = "4.30"; Volume = "  -0.67%  " }
    @{ Row = 21; Price = "326.04"; Volume = "  -0.65%  " }
    @{ Row = 22; Price = "0.999"; Volume = "  -0.24%  " }
    @{ Row = 23; Price = "5.89"; Volume = "  -1.03%  " }
    @{ Row = 24; Price = "63.62"; Volume = "  +2.35%  " }
    @{ Row = 25; Price = "0.424"; Volume = "  -4.16%  " }
    @{ Row = 26; Price = "0.169"; Volume = "  +1.73%  " }
    @{ Row = 27; Volume = "  +0.74%  " }
    @{ Row = 28; Price = "7.81"; Volume = "  -2.68%  " }
    @{ Row = 29; Coin = "PEPE"; Link = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"; Price = "0.0₃0786"; Volume = "  -1.96%  " }
    @{ Row = 30; Coin = "Aptos"; Link = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; Price = "6.75"; Volume = "  -2.74%  " }
    @{ Row = 31; Price = "1.80"; Volume = "  -0.80%  " }
    @{ Row = 32; Price = "163.71"; Volume = "  +0.72%  " }
    @{ Row = 33; Volume = "  +0.17%  " }
    @{ Row = 34; Coin = "Fetch.AI"; Link = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"; Price = "1.13"; Volume = "  -7.17%  " }
    @{ Row = 35; Coin = "ImmutableX"; Link = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; Price = "1.45"; Volume = "  -2.03%  " }
    @{ Row = 36; Price = "18.57"; Volume = "  -1.51%  " }
    @{ Row = 37; Price = "4.22"; Volume = "  -5.29%  " }
    @{ Row = 38; Price = "1.61"; Volume = "  -1.74%  " }
    @{ Row = 39; Coin = "Filecoin"; Link = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; Price = "3.69"; Volume = "  -1.24%  " }
    @{ Row = 40; Coin = "SuiNetwork"; Link = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"; Price = "0.810"; Volume = "  -3.60%  " }
    @{ Row = 41; Coin = "RenderToken"; Link = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; Price = "5.25"; Volume = "  -8.88%  " }
    @{ Row = 42; Coin = "Bittensor"; Link = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"; Price = "280.73"; Volume = "  -7.48%  " }
    @{ Row = 43; Coin = "FirstDigitalUSD"; Link = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"; Price = "0.996"; Volume = "  +0.39%  " }
    @{ Row = 44; Coin = "WhiteBITCoin"; Link = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"; Price = "10.89"; Volume = "  +0.42%  " }
    @{ Row = 45; Coin = "Mantle"; Link = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"; Price = "0.599"; Volume = "  -1.39%  " }
    @{ Row = 46; Coin = "Stellar"; Link = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; Price = "0.0938"; Volume = "  -0.02%  " }
    @{ Row = 47; Coin = "Aave"; Link = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"; Price = "125.03"; Volume = "  +0.68%  " }
    @{ Row = 48; Coin = "Hedera"; Link = "https://coinranking.com/coin/jad286TjB+hedera-hbar"; Price = "0.0514"; Volume = "  -1.01%  " }
    @{ Row = 49; Coin = "VeChain"; Link = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; Price = "0.0225"; Volume = "  -1.78%  " }
    @{ Row = 50; Coin = "InjectiveProtocol"; Link = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; Price = "17.92"; Volume = "  -2.45%  " }
    @{ Row = 51; Coin = "Maker"; Link = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"; Price = "1.776.94"; Volume = "  -2.37%  " }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("Coin"))   { $ws.Cells.Item($u.Row, 2).Value = $u.Coin }
    if ($u.ContainsKey("Link"))   { $ws.Cells.Item($u.Row, 3).Value = $u.Link }
    if ($u.ContainsKey("Price"))  {
        # Price column holds plain text (e.g. "0.999", "59.293.01"), not numbers - a
        # leading apostrophe forces Excel to keep it as text instead of auto-converting.
        $ws.Cells.Item($u.Row, 4).Value = "'" + $u.Price
    }
    if ($u.ContainsKey("Volume")) { $ws.Cells.Item($u.Row, 5).Value = $u.Volume }
}
